# NATMI LR-pair sheet update ("Natmi following Dr Hou advice"):
# the NATMI edge statistics for the Wnt2 (FAPs) -> Fzd5 pair were recomputed, and a
# new row for the "ECs" target cluster was added ahead of the existing target
# clusters (FAPs, M1, M2, Neutro, sCs), growing the sheet from 6 to 7 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster = ECs
$ws.Cells.Item(2, 1).Value = "FAPs"   # Sending cluster
$ws.Cells.Item(2, 2).Value = "Wnt2"   # Ligand symbol
$ws.Cells.Item(2, 3).Value = "Fzd5"   # Receptor symbol
$ws.Cells.Item(2, 4).Value = "ECs"   # Target cluster
$ws.Cells.Item(2, 5).Value = 3   # Ligand-expressing cells
$ws.Cells.Item(2, 6).Value = 1   # Ligand detection rate
$ws.Cells.Item(2, 7).Value = 0.574538   # Ligand average expression value
$ws.Cells.Item(2, 8).Value = 1.723614   # Ligand total expression value
$ws.Cells.Item(2, 9).Value = 1   # Ligand derived specificity of average expression value
$ws.Cells.Item(2, 10).Value = 1   # Ligand derived specificity of total expression value
$ws.Cells.Item(2, 11).Value = 1   # Receptor-expressing cells
$ws.Cells.Item(2, 12).Value = 0.5   # Receptor detection rate
$ws.Cells.Item(2, 13).Value = 2.170377   # Receptor average expression value
$ws.Cells.Item(2, 14).Value = 4.340754   # Receptor total expression value
$ws.Cells.Item(2, 15).Value = 0.1015511790371702   # Receptor derived specificity of average expression value
$ws.Cells.Item(2, 16).Value = 0.07285982038608425   # Receptor derived specificity of total expression value
$ws.Cells.Item(2, 17).Value = 1.246964060826   # Edge average expression weight
$ws.Cells.Item(2, 18).Value = 7.481784364955999   # Edge total expression weight
$ws.Cells.Item(2, 19).Value = 0.1015511790371702   # Edge average expression derived specificity
$ws.Cells.Item(2, 20).Value = 0.07285982038608425   # Edge total expression derived specificity

# Row 3: Target cluster = FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"   # Sending cluster
$ws.Cells.Item(3, 2).Value = "Wnt2"   # Ligand symbol
$ws.Cells.Item(3, 3).Value = "Fzd5"   # Receptor symbol
$ws.Cells.Item(3, 4).Value = "FAPs"   # Target cluster
$ws.Cells.Item(3, 5).Value = 3   # Ligand-expressing cells
$ws.Cells.Item(3, 6).Value = 1   # Ligand detection rate
$ws.Cells.Item(3, 7).Value = 0.574538   # Ligand average expression value
$ws.Cells.Item(3, 8).Value = 1.723614   # Ligand total expression value
$ws.Cells.Item(3, 9).Value = 1   # Ligand derived specificity of average expression value
$ws.Cells.Item(3, 10).Value = 1   # Ligand derived specificity of total expression value
$ws.Cells.Item(3, 11).Value = 3   # Receptor-expressing cells
$ws.Cells.Item(3, 12).Value = 1   # Receptor detection rate
$ws.Cells.Item(3, 13).Value = 5.061974333333333   # Receptor average expression value
$ws.Cells.Item(3, 14).Value = 15.185923   # Receptor total expression value
$ws.Cells.Item(3, 15).Value = 0.2368480046581279   # Receptor derived specificity of average expression value
$ws.Cells.Item(3, 16).Value = 0.2548966428820674   # Receptor derived specificity of total expression value
$ws.Cells.Item(3, 17).Value = 2.908296609524666   # Edge average expression weight
$ws.Cells.Item(3, 18).Value = 26.174669485722   # Edge total expression weight
$ws.Cells.Item(3, 19).Value = 0.2368480046581279   # Edge average expression derived specificity
$ws.Cells.Item(3, 20).Value = 0.2548966428820674   # Edge total expression derived specificity

# Row 4: Target cluster = M1
$ws.Cells.Item(4, 1).Value = "FAPs"   # Sending cluster
$ws.Cells.Item(4, 2).Value = "Wnt2"   # Ligand symbol
$ws.Cells.Item(4, 3).Value = "Fzd5"   # Receptor symbol
$ws.Cells.Item(4, 4).Value = "M1"   # Target cluster
$ws.Cells.Item(4, 5).Value = 3   # Ligand-expressing cells
$ws.Cells.Item(4, 6).Value = 1   # Ligand detection rate
$ws.Cells.Item(4, 7).Value = 0.574538   # Ligand average expression value
$ws.Cells.Item(4, 8).Value = 1.723614   # Ligand total expression value
$ws.Cells.Item(4, 9).Value = 1   # Ligand derived specificity of average expression value
$ws.Cells.Item(4, 10).Value = 1   # Ligand derived specificity of total expression value
$ws.Cells.Item(4, 11).Value = 3   # Receptor-expressing cells
$ws.Cells.Item(4, 12).Value = 1   # Receptor detection rate
$ws.Cells.Item(4, 13).Value = 3.815520666666667   # Receptor average expression value
$ws.Cells.Item(4, 14).Value = 11.446562   # Receptor total expression value
$ws.Cells.Item(4, 15).Value = 0.1785268745202745   # Receptor derived specificity of average expression value
$ws.Cells.Item(4, 16).Value = 0.1921312406457904   # Receptor derived specificity of total expression value
$ws.Cells.Item(4, 17).Value = 2.192161612785333   # Edge average expression weight
$ws.Cells.Item(4, 18).Value = 19.729454515068   # Edge total expression weight
$ws.Cells.Item(4, 19).Value = 0.1785268745202745   # Edge average expression derived specificity
$ws.Cells.Item(4, 20).Value = 0.1921312406457904   # Edge total expression derived specificity

# Row 5: Target cluster = M2
$ws.Cells.Item(5, 1).Value = "FAPs"   # Sending cluster
$ws.Cells.Item(5, 2).Value = "Wnt2"   # Ligand symbol
$ws.Cells.Item(5, 3).Value = "Fzd5"   # Receptor symbol
$ws.Cells.Item(5, 4).Value = "M2"   # Target cluster
$ws.Cells.Item(5, 5).Value = 3   # Ligand-expressing cells
$ws.Cells.Item(5, 6).Value = 1   # Ligand detection rate
$ws.Cells.Item(5, 7).Value = 0.574538   # Ligand average expression value
$ws.Cells.Item(5, 8).Value = 1.723614   # Ligand total expression value
$ws.Cells.Item(5, 9).Value = 1   # Ligand derived specificity of average expression value
$ws.Cells.Item(5, 10).Value = 1   # Ligand derived specificity of total expression value
$ws.Cells.Item(5, 11).Value = 3   # Receptor-expressing cells
$ws.Cells.Item(5, 12).Value = 1   # Receptor detection rate
$ws.Cells.Item(5, 13).Value = 4.666218666666667   # Receptor average expression value
$ws.Cells.Item(5, 14).Value = 13.998656   # Receptor total expression value
$ws.Cells.Item(5, 15).Value = 0.2183307357409577   # Receptor derived specificity of average expression value
$ws.Cells.Item(5, 16).Value = 0.2349682939430755   # Receptor derived specificity of total expression value
$ws.Cells.Item(5, 17).Value = 2.680919940309333   # Edge average expression weight
$ws.Cells.Item(5, 18).Value = 24.128279462784   # Edge total expression weight
$ws.Cells.Item(5, 19).Value = 0.2183307357409577   # Edge average expression derived specificity
$ws.Cells.Item(5, 20).Value = 0.2349682939430755   # Edge total expression derived specificity

# Row 6: Target cluster = Neutro
$ws.Cells.Item(6, 1).Value = "FAPs"   # Sending cluster
$ws.Cells.Item(6, 2).Value = "Wnt2"   # Ligand symbol
$ws.Cells.Item(6, 3).Value = "Fzd5"   # Receptor symbol
$ws.Cells.Item(6, 4).Value = "Neutro"   # Target cluster
$ws.Cells.Item(6, 5).Value = 3   # Ligand-expressing cells
$ws.Cells.Item(6, 6).Value = 1   # Ligand detection rate
$ws.Cells.Item(6, 7).Value = 0.574538   # Ligand average expression value
$ws.Cells.Item(6, 8).Value = 1.723614   # Ligand total expression value
$ws.Cells.Item(6, 9).Value = 1   # Ligand derived specificity of average expression value
$ws.Cells.Item(6, 10).Value = 1   # Ligand derived specificity of total expression value
$ws.Cells.Item(6, 11).Value = 3   # Receptor-expressing cells
$ws.Cells.Item(6, 12).Value = 1   # Receptor detection rate
$ws.Cells.Item(6, 13).Value = 3.288577   # Receptor average expression value
$ws.Cells.Item(6, 14).Value = 9.865731   # Receptor total expression value
$ws.Cells.Item(6, 15).Value = 0.1538713650690733   # Receptor derived specificity of average expression value
$ws.Cells.Item(6, 16).Value = 0.1655968959856798   # Receptor derived specificity of total expression value
$ws.Cells.Item(6, 17).Value = 1.889412452426   # Edge average expression weight
$ws.Cells.Item(6, 18).Value = 17.004712071834   # Edge total expression weight
$ws.Cells.Item(6, 19).Value = 0.1538713650690733   # Edge average expression derived specificity
$ws.Cells.Item(6, 20).Value = 0.1655968959856798   # Edge total expression derived specificity

# Row 7: Target cluster = sCs
$ws.Cells.Item(7, 1).Value = "FAPs"   # Sending cluster
$ws.Cells.Item(7, 2).Value = "Wnt2"   # Ligand symbol
$ws.Cells.Item(7, 3).Value = "Fzd5"   # Receptor symbol
$ws.Cells.Item(7, 4).Value = "sCs"   # Target cluster
$ws.Cells.Item(7, 5).Value = 3   # Ligand-expressing cells
$ws.Cells.Item(7, 6).Value = 1   # Ligand detection rate
$ws.Cells.Item(7, 7).Value = 0.574538   # Ligand average expression value
$ws.Cells.Item(7, 8).Value = 1.723614   # Ligand total expression value
$ws.Cells.Item(7, 9).Value = 1   # Ligand derived specificity of average expression value
$ws.Cells.Item(7, 10).Value = 1   # Ligand derived specificity of total expression value
$ws.Cells.Item(7, 11).Value = 2   # Receptor-expressing cells
$ws.Cells.Item(7, 12).Value = 1   # Receptor detection rate
$ws.Cells.Item(7, 13).Value = 2.3695805   # Receptor average expression value
$ws.Cells.Item(7, 14).Value = 4.739161   # Receptor total expression value
$ws.Cells.Item(7, 15).Value = 0.1108718409743963   # Receptor derived specificity of average expression value
$ws.Cells.Item(7, 16).Value = 0.07954710615730251   # Receptor derived specificity of total expression value
$ws.Cells.Item(7, 17).Value = 1.361414041309   # Edge average expression weight
$ws.Cells.Item(7, 18).Value = 8.168484247854   # Edge total expression weight
$ws.Cells.Item(7, 19).Value = 0.1108718409743963   # Edge average expression derived specificity
$ws.Cells.Item(7, 20).Value = 0.07954710615730251   # Edge total expression derived specificity
